$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of kaspa buy data appended after the run on 2025-10-10.
# The Date column in this sheet is stored as plain text (e.g. "08/29/2025"),
# so force text formatting on A9 before assigning the value to stop Excel
# from auto-converting the "10/10/2025" string into a date serial number.
$dateCell = $ws.Cells.Item(9, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "10/10/2025"
# Drop the temporary Text number-format style again so the new cell ends up
# with the same (default) style as its siblings in column A.
$dateCell.Style = "Normal"

$ws.Cells.Item(9, 2).Value = 335.1719999999996
$ws.Cells.Item(9, 3).Value = 0.1491771389018178
$ws.Cells.Item(9, 4).Value = 25
